# Messenger_Thesis_1._0/wwwroot/Book1.xlsx edit
# Replaces the placeholder "Account name/Address/Area" table (9 data rows)
# with a new "NAMe/ADDRESS/AREA" table of names and cities (7 rows incl.
# header), resizes column A, and updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table shrinks from 9 data rows to 7 rows (header + 6 rows) -- drop
# the two trailing rows first.
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(8).Delete()

# Header row
$ws.Range("A1").Value = "NAMe"
$ws.Range("B1").Value = "ADDRESS"
$ws.Range("C1").Value = "AREA"

# Row 2 entered directly
$ws.Cells.Item(2, 1).Value = "Alberto Bugatsing"
$ws.Cells.Item(2, 2).Value = "Manila"
$ws.Cells.Item(2, 3).Value = "Manila"

# Rest of column A (rows 3-7)
$names = @("Areeya Mae", "Jinkee Pacquiao", "Shania Twain", "Taylor Swift", "Mitchel B.")
$row = 3
foreach ($n in $names) {
    $ws.Cells.Item($row, 1).Value = $n
    $row++
}

# Rest of columns B & C (rows 3-7)
$addrs = @("Caloocan", "Anna", "Pasig", "Pasay", "Mandaluyong")
$row = 3
foreach ($a in $addrs) {
    $ws.Cells.Item($row, 2).Value = $a
    $ws.Cells.Item($row, 3).Value = $a
    $row++
}

# Widen column A to fit the longer names
$ws.Columns.Item(1).ColumnWidth = 40.6

# Active selection moves to G5
$ws.Range("G5").Select() | Out-Null
